$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.678.31"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "2.292.66"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'96.36"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "'268.54"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.606"
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").Value = "'45.22"
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").Value = "'0.0933"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "'7.87"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "2.636.78"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "'0.852"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "2.302.08"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "43.668.94"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'0.0000113"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "'71.98"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("D22").Value = "'2.52"
$ws.Range("E22").Value = "  +10.92%  "
$ws.Range("D23").Value = "'232.28"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("D24").Value = "'9.08"
$ws.Range("E24").Value = "  -6.03%  "
$ws.Range("E25").Value = "  +5.52%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'11.22"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D30").Value = "'38.43"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'174.62"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").Value = "'21.80"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'0.0901"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'4.51"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "'0.238"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'12.11"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'64.65"
$ws.Range("E43").Value = "  +4.09%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'1.34"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").Value = "'8.77"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "'5.15"
$ws.Range("E47").Value = "  -5.73%  "
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "'97.20"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.54"
$ws.Range("E50").Value = "  +12.65%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.441"
$ws.Range("E51").Value = "  +4.43%  "
